# AnimationInfo.xlsx update: a new "Idle02" animation clip is inserted
# right after "Idle01" in the table, so "Wait" and "Walk" each move down
# one row and every row gets new Start/End frame ranges.
#
# Final table (Sheet1, columns C:F, header on row 3):
#   row4: Index=0  Name=Idle01  Start=0    End=128
#   row5: Index=1  Name=Idle02  Start=129  End=226   (new row)
#   row6: Index=2  Name=Wait    Start=227  End=347
#   row7: Index=3  Name=Walk    Start=348  End=378   (new row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row right after the existing last data row (row 6). This
# both makes room for the "Walk" row that now comes last, and lets the
# new row naturally pick up the same formatting/style as the rows above
# it instead of us having to copy a style over explicitly.
$ws.Rows("7").Insert()

# Row 4: was Animation Index 0 = Wait -> now Idle01
$ws.Range("D4").Value = "Idle01"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 128

# Row 5: was Animation Index 1 = Idle01 -> now the brand-new Idle02 entry
$ws.Range("D5").Value = "Idle02"
$ws.Range("E5").Value = 129
$ws.Range("F5").Value = 226

# Row 6: was Animation Index 2 = Walk -> now Wait
$ws.Range("D6").Value = "Wait"
$ws.Range("E6").Value = 227
$ws.Range("F6").Value = 347

# Row 7 (new): Animation Index 3 = Walk, moved down from row 6
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Walk"
$ws.Range("E7").Value = 348
$ws.Range("F7").Value = 378

# Leave the selection on F7, matching the author's final cursor position
$ws.Range("F7").Select()
